$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4 data, written in the same column order the original author
# used (C, A, B, D, E) so the shared-string table indices line up with
# the target workbook.
$ws.Range("C4").Value = ' It isssss rare to come by\nnewssssss on [CS:P]Zero Isle[CR]. Right, bossssssss?'
$ws.Range("A4").Value = 'SCRIPT/G01P03A/um1105.ssb '
$ws.Range("B4").Value = 63
$ws.Range("D4").Value = ' Ссссс [CS:P]Нуль-Оссссстровов[CR] нет\nпочти никаких вессссстей. Верно, боссссс?'
$ws.Range("E4").Value = ' Òòòòò [CS:P]Îôìû-Ïòòòòòóñïâïâ[CR] îåó\nðïœóé îéëàëéö âåòòòòòóåê. Âåñîï, áïòòòòò?'

# Match the style used by row 3 (font + wrap, bottom-thin border) as a
# base, then give row 4 its own top+bottom thin border since it is now
# the last row of the table.
$ws.Range("A3:B3").Copy()
$ws.Range("A4:B4").PasteSpecial(-4122)
$ws.Range("C3:E3").Copy()
$ws.Range("C4:E4").PasteSpecial(-4122)

$ws.Range("A4:E4").Borders.Item(8).LineStyle = 1
$ws.Range("A4:E4").Borders.Item(8).Weight = 2

$ws.Rows.Item(4).RowHeight = 43.2

$ws.Range("A4:E4").Select()
